$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Image Caption) so the new, longer captions are readable.
# The OOXML <col> width stores a value 0.8333333333333333 (5/6) higher than
# the ColumnWidth COM property, so subtract that offset to land on width=82.
$ws.Columns.Item(5).ColumnWidth = 82 - 5/6

# For each row, update the Image Caption (E), and where applicable the
# Caption Sentiment (F), Overall Sentiment (G) and Confidence (H) columns
# with the freshly generated captions/scores. Cells that already hold the
# right value in the source workbook are left untouched.
# The Confidence column stores percentages as literal text (e.g. "50.0%"),
# so NumberFormat is forced to "@" (Text) before assignment to stop Excel
# from auto-converting the string into a numeric percentage.

# Row 2
$ws.Range("E2").Value = "a man in a white shirt and black pants stands at a podium with a microphone"
$ws.Range("F2").Value = "Neutral"
$ws.Range("G2").Value = "Neutral"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "50.0%"

# Row 3
$ws.Range("E3").Value = "a group of people posing for a photo"
$ws.Range("F3").Value = "Neutral"
$ws.Range("G3").Value = "Neutral"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "50.0%"

# Row 4
$ws.Range("E4").Value = "a man in a yellow shirt is holding a microphone"
$ws.Range("F4").Value = "Neutral"
$ws.Range("G4").Value = "Positive"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "68.06%"

# Row 5
$ws.Range("E5").Value = "obama and dute"
$ws.Range("F5").Value = "Neutral"
$ws.Range("G5").Value = "Neutral"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "50.0%"

# Row 6
$ws.Range("E6").Value = "person giving a speech on stage"
$ws.Range("F6").Value = "Positive"
$ws.Range("G6").Value = "Positive"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "67.0%"

# Row 7
$ws.Range("E7").Value = "a woman in a red jacket and glasses speaking into a microphone"
$ws.Range("F7").Value = "Neutral"
$ws.Range("G7").Value = "Neutral"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "50.0%"

# Row 8
$ws.Range("E8").Value = "a man in a white shirt and a hat with his head in his hands, with the words what"
$ws.Range("F8").Value = "Neutral"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "100%"

# Row 9
$ws.Range("E9").Value = "a man in a white shirt is crying and holding his head"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "73.84%"

# Row 10
$ws.Range("E10").Value = "two men in suits standing next to each other men"
$ws.Range("F10").Value = "Neutral"

# Row 11
$ws.Range("E11").Value = "a group of people riding in a boat down a flooded street"
$ws.Range("F11").Value = "Neutral"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "74.97%"
